$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A13 holds a date-like string ("2022-03-24" -> "2022-02-23"). The cell keeps
# its pre-existing date NumberFormat but the underlying value is stored as
# plain text in the workbook, so a direct Value assignment must be avoided
# (Excel silently reinterprets a date-looking string as a real date serial
# when assigned straight through .Value). Using a text formula and then
# collapsing it to a static value via copy/paste-special keeps it as text
# without touching any cell styles.
$ws.Range("A13").Formula = '="2022-02-23"'
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A15").Value = "Glee Star Enterprises"
$ws.Range("A16").Value = "Glee Star Enterprises`n"
$ws.Range("A17").Value = "102 Centennial II Extension St., Pinagbuhatan, Pasig City`n"
$ws.Range("A20").Value = "Dear Glee Star Enterprises`n"
$ws.Range("A22").Value = "We are pleased to inform you that your Quotation for the Procurement of  for the LOREM IPSUMwith  Purchase Order equivalent to Php 147,200.00 is hereby accepted. "
$ws.Range("A41").Value = "                                  Glee Star Enterprises"
